# Implement ICloneable on LinearRegressionResult
# Updates the "ICloneable" column (F) for LinearRegressionResult (row 6) and
# DecimalLinearRegressionResult (row 7) from "TODO" to "Oui", reusing the
# existing "Good" (Bueno) cell style already used by sibling cells such as B6/B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Update F6 (LinearRegressionResult / ICloneable) ---
$ws.Range("B6").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F6").Value = "Oui"

# --- Update F7 (DecimalLinearRegressionResult / ICloneable) ---
$ws.Range("B7").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F7").Value = "Oui"

$excel.CutCopyMode = 0

# --- Reset zoom to 100% (was 130%) ---
$excel.ActiveWindow.Zoom = 100

# --- Move the active selection on the frozen (bottom-right) pane to F10 ---
$ws.Range("F10").Select()
